$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text while we update values, to avoid
# Excel auto-converting numeric-looking strings (e.g. "1.003", "309.06") into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.097.27"
$ws.Range("E2").Value = "  -2.45%  "

$ws.Range("D3").Value = "1.638.87"
$ws.Range("E3").Value = "  -2.39%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "309.06"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "0.3935"
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("D8").Value = "0.3857"
$ws.Range("E8").Value = "  -2.45%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "50.13"
$ws.Range("E10").Value = "  -3.77%  "

$ws.Range("D11").Value = "1.369"
$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("D12").Value = "0.08552"
$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("E13").Value = "  -6.01%  "

$ws.Range("D14").Value = "7.068"
$ws.Range("E14").Value = "  -3.41%  "

$ws.Range("D15").Value = "0.00001282"
$ws.Range("E15").Value = "  -2.53%  "

$ws.Range("D16").Value = "7.494"
$ws.Range("E16").Value = "  -3.45%  "

$ws.Range("D17").Value = "1.636.58"
$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").Value = "93.71"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").Value = "0.06919"
$ws.Range("E19").Value = "  -1.97%  "

$ws.Range("D20").Value = "20.30"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").Value = "6.915"
$ws.Range("E21").Value = "  -2.19%  "

$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("D24").Value = "24.112.55"
$ws.Range("E24").Value = "  -2.38%  "

$ws.Range("D25").Value = "2.409"

$ws.Range("D26").Value = "2.875"
$ws.Range("E26").Value = "  +3.16%  "

$ws.Range("D27").Value = "22.19"
$ws.Range("E27").Value = "  -5.07%  "

$ws.Range("D28").Value = "158.24"
$ws.Range("E28").Value = "  -2.48%  "

$ws.Range("D29").Value = "139.78"
$ws.Range("E29").Value = "  -4.98%  "

$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "5.269"
$ws.Range("E30").Value = "  -9.81%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "7.993"
$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("D32").Value = "2.480"
$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("D33").Value = "1.822.13"
$ws.Range("E33").Value = "  -2.29%  "

$ws.Range("D34").Value = "0.08058"
$ws.Range("E34").Value = "  -4.14%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "6.730"
$ws.Range("E35").Value = "  -2.82%  "

$ws.Range("D36").Value = "0.02903"
$ws.Range("E36").Value = "  -4.63%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.9688"
$ws.Range("E37").Value = "  -2.89%  "

$ws.Range("D38").Value = "0.2689"
$ws.Range("E38").Value = "  -4.49%  "

$ws.Range("D39").Value = "0.09216"
$ws.Range("E39").Value = "  -2.83%  "

$ws.Range("D40").Value = "10.34"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("D41").Value = "1.428"
$ws.Range("E41").Value = "  -5.06%  "

$ws.Range("D42").Value = "0.7517"
$ws.Range("E42").Value = "  -5.28%  "

$ws.Range("D43").Value = "13.07"
$ws.Range("E43").Value = "  -3.68%  "

$ws.Range("D44").Value = "16.18"
$ws.Range("E44").Value = "  -2.94%  "

$ws.Range("D45").Value = "0.6906"
$ws.Range("E45").Value = "  -3.25%  "

$ws.Range("D46").Value = "2.456"
$ws.Range("E46").Value = "  -4.21%  "

$ws.Range("D47").Value = "4.090"
$ws.Range("E47").Value = "  -2.60%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").Value = "0.08333"
$ws.Range("E49").Value = "  -4.13%  "

$ws.Range("D50").Value = "1.264"
$ws.Range("E50").Value = "  -6.12%  "

$ws.Range("D51").Value = "133.43"
$ws.Range("E51").Value = "  -3.32%  "

# Restore default (General/Normal) style on column D so formatting matches the original workbook.
$priceRange.Style = "Normal"
